# Update "paises" workbook (Pais sheet):
#  - Re-sort three countries around Tayikistan/Haiti/Libia: Libia overtakes
#    Tayikistan and Haiti in total cases, so its row (A95/B95/.../H95) now
#    shows Libia's updated numbers, while Tayikistan and Haiti shift down
#    into rows 96 and 97 respectively (values cascade down one row).
#  - Swap Montserrat ahead of Islas Malvinas (rows 213/214) - same total
#    cases (13), but Montserrat's breakdown now sorts first.
#  - Refresh a handful of per-country case/recovered/death counts.
#  - Update the "Datos actualizados..." timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 04:59"

# --- Row 29: Kazajistan ---
$ws.Range("B29").Value = 103033
$ws.Range("C29").Value = 337
$ws.Range("D29").Value = 81558
$ws.Range("E29").Value = 20206

# --- Row 31: Bolivia ---
$ws.Range("B31").Value = 100344
$ws.Range("C31").Value = 1198
$ws.Range("D31").Value = 36491
$ws.Range("E31").Value = 59795
$ws.Range("G31").Value = 55
$ws.Range("H31").Value = 4058

# --- Row 40: Belgica ---
$ws.Range("B40").Value = 78323
$ws.Range("C40").Value = 454
$ws.Range("D40").Value = 17994
$ws.Range("E40").Value = 50390
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 9939

# --- Row 72: Australia ---
$ws.Range("D72").Value = 14080
$ws.Range("E72").Value = 8811

# --- Rows 95-97: Libia overtakes Tayikistan & Haiti (cascading re-sort) ---
# Row 95 becomes Libia with its fresh numbers.
$ws.Range("A95").Value = "Libia"
$ws.Range("B95").Value = 8172
$ws.Range("D95").Value = 933
$ws.Range("E95").Value = 7086
$ws.Range("H95").Value = 153

# Row 96 becomes Tayikistan, carrying the old row-95 numbers.
$ws.Range("A96").Value = "Tayikistan"
$ws.Range("B96").Value = 8065
$ws.Range("D96").Value = 6855
$ws.Range("E96").Value = 1146
$ws.Range("H96").Value = 64

# Row 97 becomes Haiti, carrying the old row-96 numbers.
$ws.Range("A97").Value = "Haiti"
$ws.Range("B97").Value = 7879
$ws.Range("D97").Value = 5235
$ws.Range("E97").Value = 2448
$ws.Range("H97").Value = 196

# --- Row 141: Nueva Zelanda ---
$ws.Range("B141").Value = 1631
$ws.Range("C141").Value = 9
$ws.Range("E141").Value = 78

# --- Row 158: Vietnam ---
$ws.Range("B158").Value = 964
$ws.Range("C158").Value = 2
$ws.Range("E158").Value = 484

# --- Row 179: San Martin (Parte Holandesa) ---
$ws.Range("B179").Value = 317
$ws.Range("C179").Value = 17
$ws.Range("E179").Value = 193

# --- Row 194: Antigua y Barbuda ---
$ws.Range("D194").Value = 88
$ws.Range("E194").Value = 2

# --- Rows 213-214: Montserrat now sorts ahead of Islas Malvinas ---
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
